$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129, shifting existing rows 129-131 down to 130-132
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new weekly data entry,
# carrying over the static descriptive columns from the row that used to occupy
# row 129 (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Variedad,
# Calidad, Unidad de comercializacion, Origen, Kg o Unidades, Clasificacion)
# and updating the date/volume/price columns with the new values.
$ws.Cells.Item(129, 1).Value = 10
$ws.Cells.Item(129, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(129, 3).Value = "La Araucanía"
$ws.Cells.Item(129, 4).Value = 44656
$ws.Cells.Item(129, 5).Value = 9
$ws.Cells.Item(129, 6).Value = 100112012
$ws.Cells.Item(129, 7).Value = "Espinaca"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 95
$ws.Cells.Item(129, 11).Value = 9000
$ws.Cells.Item(129, 12).Value = 9000
$ws.Cells.Item(129, 13).Value = 9000
$ws.Cells.Item(129, 14).Value = "$/docena de atados"
$ws.Cells.Item(129, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(129, 16).Value = 3000
$ws.Cells.Item(129, 17).Value = 3
$ws.Cells.Item(129, 18).Value = "Hortaliza"
